# update diagrams and images for clarity and consistency
#
# The "Refinitiv Real-Time" label shows up twice on slide 1 (two
# rectangle shapes in the architecture diagram) and both need to be
# renamed to "Real-Time Platform".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange

        if ($tr.Text -eq "Refinitiv Real-Time") {
            $tr.Text = "Real-Time Platform"
        }
    }
}
